# Insert a new weekly record at row 55 (Albahaca, Feria Lagunitas de Puerto Montt).
# This pushes the previous rows 55-168 down to 56-169 (Excel's native Insert behavior),
# and the freshly inserted row 55 is populated with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 55:168 down to 56:169, creating a new blank row 55.
$ws.Rows("55:55").Insert()

# Populate the newly inserted row 55 with the new weekly record.
$ws.Range("A55").Value = 4
$ws.Range("B55").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C55").Value = "Los Lagos"
$ws.Range("D55").Value = 44967
$ws.Range("E55").Value = 10
$ws.Range("F55").Value = 100112052
$ws.Range("G55").Value = "Albahaca"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 90
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 6000
$ws.Range("M55").Value = 6000
$ws.Range("N55").Value = "$/docena de matas"
$ws.Range("O55").Value = "Región Metropolitana"
$ws.Range("P55").Value = 1000
$ws.Range("Q55").Value = 6
$ws.Range("R55").Value = "Hortaliza"

# Preserve the date number format on the new row's date cell.
$ws.Range("D55").NumberFormat = $ws.Range("D56").NumberFormat
